$d = $word.ActiveDocument

# Adds a brand-new empty paragraph at the very end of the document body
# (collapsing Document.Content to its end and assigning a bare paragraph
# mark keeps the serialized paragraph free of any stray empty run).
function Add-EmptyParagraph {
    $endRange = $d.Content
    $endRange.Collapse(0)
    $endRange.Text = "`r"
}

# Appends a new paragraph holding $text right after the current last
# paragraph of the document.
function Add-TextParagraph([string]$text) {
    $lastRange = $d.Paragraphs.Last.Range
    $lastRange.Collapse(0)
    $lastRange.InsertAfter("`r" + $text)
}

Add-EmptyParagraph
Add-TextParagraph "Incluir variables año niño o niña según el preddición si será un año con dicha particualaridad climática."
Add-EmptyParagraph
Add-TextParagraph "Incluir variable del tiempo de cosecha del producto"
Add-TextParagraph "Temperatura"
Add-TextParagraph "Velocidad del viento"
Add-TextParagraph "Temporadas de cultivo por tipo de cultivo"
Add-EmptyParagraph
Add-TextParagraph "Determinar las coberturas a partir de las variables más relevantes. Ejemplo si la variable velocidad del viento es representativa identifi carla como una cobertura."
Add-EmptyParagraph
